$wb = $excel.ActiveWorkbook

# --- Update localization status text ("Ready for handoff" -> "In Translation") ---
# This status string shows up once per language column on the "Overview" sheet
# (one cell per language) and once per language-specific sheet in its "Status"
# column (row 2, the only data row in this report).

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the per-language status columns to fit the shorter label ---
# Columns E & F on "Overview" and column C on each language sheet all hold the
# same status value, and all shrink from the old wider width to a narrower one.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
